# Table styles and scroll
# Refactored styles for Directory table with using anothers styles.
# A bit hardcoding but with better visual interface.
#
# - Re-keyed a handful of vehicle plate numbers in the "Довідники" sheet.
# - Renamed the subdivision label (column K) used by the existing rows.
# - Filled in the previously-empty 9th directory row with a new vehicle
#   entry, copying formatting/types from the row above so the new cells
#   stay text (matching the rest of the table) instead of turning numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Довідники")

# --- Update existing license-plate values -------------------------------
$ws.Range("B3").Value = "BB1234ЗA"
$ws.Range("B4").Value = "BB1234ЗX"
$ws.Range("B5").Value = "BB12342X"
$ws.Range("B6").Value = "BB12344X"
$ws.Range("B7").Value = "BB12345X"
$ws.Range("B8").Value = "BB12346X"

# --- Rename subdivision label in column K for all existing rows --------
$ws.Range("K3:K8").Value = "А0000 (ПВЗ)"

# --- Populate the previously-empty row 9 with a new directory entry ----
# Copy row 8's formatting/types down into row 9 first so the text-like
# numeric columns (D, F) keep being stored as text, then overwrite the
# cells that actually differ for the new entry.
$ws.Range("A8:M8").Copy()
$ws.Range("A9:M9").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("B9").Value = "BB12347X"
